$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Gnai2"
$ws.Range("C2").Value = "Tbxa2r"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.999957666666667
$ws.Range("N2").Value = 8.999873000000001
$ws.Range("O2").Value = 0.5075658552021639
$ws.Range("P2").Value = 0.5075658552021638
$ws.Range("Q2").Value = 604.3107003842273
$ws.Range("R2").Value = 5438.796303458045
$ws.Range("S2").Value = 0.2453319679426353
$ws.Range("T2").Value = 0.2453319679426353

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Gnai2"
$ws.Range("C3").Value = "Tbxa2r"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.303808666666667
$ws.Range("N3").Value = 6.911426000000001
$ws.Range("O3").Value = 0.389783705654121
$ws.Range("P3").Value = 0.3897837056541209
$ws.Range("Q3").Value = 464.0786249665699
$ws.Range("R3").Value = 4176.707624699129
$ws.Range("S3").Value = 0.1884019632132471
$ws.Range("T3").Value = 0.1884019632132471

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Gnai2"
$ws.Range("C4").Value = "Tbxa2r"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.448687
$ws.Range("N4").Value = 1.346061
$ws.Range("O4").Value = 0.07591380485249957
$ws.Range("P4").Value = 0.07591380485249956
$ws.Range("Q4").Value = 90.38339381787868
$ws.Range("R4").Value = 813.4505443609081
$ws.Range("S4").Value = 0.03669293934490316
$ws.Range("T4").Value = 0.03669293934490316

$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Gnai2"
$ws.Range("C5").Value = "Tbxa2r"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1580263333333334
$ws.Range("N5").Value = 0.474079
$ws.Range("O5").Value = 0.02673663429121574
$ws.Range("P5").Value = 0.02673663429121573
$ws.Range("Q5").Value = 31.83278392122356
$ws.Range("R5").Value = 286.495055291012
$ws.Range("S5").Value = 0.01292315280785369
$ws.Range("T5").Value = 0.01292315280785369

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Gnai2"
$ws.Range("C6").Value = "Tbxa2r"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 65.41736466666667
$ws.Range("H6").Value = 196.252094
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.999957666666667
$ws.Range("N6").Value = 8.999873000000001
$ws.Range("O6").Value = 0.5075658552021639
$ws.Range("P6").Value = 0.5075658552021638
$ws.Range("Q6").Value = 196.2493246648958
$ws.Range("R6").Value = 1766.243921984062
$ws.Range("S6").Value = 0.07967132304101211
$ws.Range("T6").Value = 0.07967132304101211

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Gnai2"
$ws.Range("C7").Value = "Tbxa2r"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 65.41736466666667
$ws.Range("H7").Value = 196.252094
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.303808666666667
$ws.Range("N7").Value = 6.911426000000001
$ws.Range("O7").Value = 0.389783705654121
$ws.Range("P7").Value = 0.3897837056541209
$ws.Range("Q7").Value = 150.7090916695605
$ws.Range("R7").Value = 1356.381825026044
$ws.Range("S7").Value = 0.06118335820072685
$ws.Range("T7").Value = 0.06118335820072685

$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Gnai2"
$ws.Range("C8").Value = "Tbxa2r"
$ws.Range("D8").Value = "MuSCs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 65.41736466666667
$ws.Range("H8").Value = 196.252094
$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.448687
$ws.Range("N8").Value = 1.346061
$ws.Range("O8").Value = 0.07591380485249957
$ws.Range("P8").Value = 0.07591380485249956
$ws.Range("Q8").Value = 29.35192110019267
$ws.Range("R8").Value = 264.167289901734
$ws.Range("S8").Value = 0.01191599712172692
$ws.Range("T8").Value = 0.01191599712172692

$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Gnai2"
$ws.Range("C9").Value = "Tbxa2r"
$ws.Range("D9").Value = "Resolving-Mac"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 65.41736466666667
$ws.Range("H9").Value = 196.252094
$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1580263333333334
$ws.Range("N9").Value = 0.474079
$ws.Range("O9").Value = 0.02673663429121574
$ws.Range("P9").Value = 0.02673663429121573
$ws.Range("Q9").Value = 10.33766627460289
$ws.Range("R9").Value = 93.038996471426
$ws.Range("S9").Value = 0.00419678157191329
$ws.Range("T9").Value = 0.00419678157191329

$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Gnai2"
$ws.Range("C10").Value = "Tbxa2r"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.999957666666667
$ws.Range("N10").Value = 8.999873000000001
$ws.Range("O10").Value = 0.5075658552021639
$ws.Range("P10").Value = 0.5075658552021638
$ws.Range("Q10").Value = 181.3019655917169
$ws.Range("R10").Value = 1631.717690325452
$ws.Range("S10").Value = 0.0736031448428822
$ws.Range("T10").Value = 0.0736031448428822

$ws.Range("A11").Value = "MuSCs"
$ws.Range("B11").Value = "Gnai2"
$ws.Range("C11").Value = "Tbxa2r"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 2.303808666666667
$ws.Range("N11").Value = 6.911426000000001
$ws.Range("O11").Value = 0.389783705654121
$ws.Range("P11").Value = 0.3897837056541209
$ws.Range("Q11").Value = 139.2303112323582
$ws.Range("R11").Value = 1253.072801091224
$ws.Range("S11").Value = 0.05652331860114714
$ws.Range("T11").Value = 0.05652331860114714

$ws.Range("A12").Value = "MuSCs"
$ws.Range("B12").Value = "Gnai2"
$ws.Range("C12").Value = "Tbxa2r"
$ws.Range("D12").Value = "MuSCs"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.448687
$ws.Range("N12").Value = 1.346061
$ws.Range("O12").Value = 0.07591380485249957
$ws.Range("P12").Value = 0.07591380485249956
$ws.Range("Q12").Value = 27.11632765332934
$ws.Range("R12").Value = 244.046948879964
$ws.Range("S12").Value = 0.01100841342431775
$ws.Range("T12").Value = 0.01100841342431775

$ws.Range("A13").Value = "MuSCs"
$ws.Range("B13").Value = "Gnai2"
$ws.Range("C13").Value = "Tbxa2r"
$ws.Range("D13").Value = "Resolving-Mac"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("K13").Value = 2
$ws.Range("L13").Value = 0.6666666666666666
$ws.Range("M13").Value = 0.1580263333333334
$ws.Range("N13").Value = 0.474079
$ws.Range("O13").Value = 0.02673663429121574
$ws.Range("P13").Value = 0.02673663429121573
$ws.Range("Q13").Value = 9.550296381488446
$ws.Range("R13").Value = 85.95266743339602
$ws.Range("S13").Value = 0.003877133077763292
$ws.Range("T13").Value = 0.003877133077763292

$ws.Range("A14").Value = "Resolving-Mac"
$ws.Range("B14").Value = "Gnai2"
$ws.Range("C14").Value = "Tbxa2r"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 2.999957666666667
$ws.Range("N14").Value = 8.999873000000001
$ws.Range("O14").Value = 0.5075658552021639
$ws.Range("P14").Value = 0.5075658552021638
$ws.Range("Q14").Value = 268.3928376254028
$ws.Range("R14").Value = 2415.535538628626
$ws.Range("S14").Value = 0.1089594193756342
$ws.Range("T14").Value = 0.1089594193756342

$ws.Range("A15").Value = "Resolving-Mac"
$ws.Range("B15").Value = "Gnai2"
$ws.Range("C15").Value = "Tbxa2r"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 2.303808666666667
$ws.Range("N15").Value = 6.911426000000001
$ws.Range("O15").Value = 0.389783705654121
$ws.Range("P15").Value = 0.3897837056541209
$ws.Range("Q15").Value = 206.1114902596945
$ws.Range("R15").Value = 1855.00341233725
$ws.Range("S15").Value = 0.0836750656389998
$ws.Range("T15").Value = 0.0836750656389998

$ws.Range("A16").Value = "Resolving-Mac"
$ws.Range("B16").Value = "Gnai2"
$ws.Range("C16").Value = "Tbxa2r"
$ws.Range("D16").Value = "MuSCs"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.448687
$ws.Range("N16").Value = 1.346061
$ws.Range("O16").Value = 0.07591380485249957
$ws.Range("P16").Value = 0.07591380485249956
$ws.Range("Q16").Value = 40.14202549379167
$ws.Range("R16").Value = 361.278229444125
$ws.Range("S16").Value = 0.01629645496155174
$ws.Range("T16").Value = 0.01629645496155174

$ws.Range("A17").Value = "Resolving-Mac"
$ws.Range("B17").Value = "Gnai2"
$ws.Range("C17").Value = "Tbxa2r"
$ws.Range("D17").Value = "Resolving-Mac"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("K17").Value = 2
$ws.Range("L17").Value = 0.6666666666666666
$ws.Range("M17").Value = 0.1580263333333334
$ws.Range("N17").Value = 0.474079
$ws.Range("O17").Value = 0.02673663429121574
$ws.Range("P17").Value = 0.02673663429121573
$ws.Range("Q17").Value = 14.13791150926389
$ws.Range("R17").Value = 127.241203583375
$ws.Range("S17").Value = 0.005739566833685464
$ws.Range("T17").Value = 0.005739566833685463
